# Bootstrap 4.0 "novedades" slides — Dropdowns section:
# "Ya no es necesario añadir [ ][.]" (three separate runs) becomes a single
# run whose text embeds the literal HTML snippet
# `<span class="caret"></span>` right before the trailing period.

$d = $word.ActiveDocument

# Locate the sentence (it currently reads as "Ya no es necesario añadir .",
# the space and the final period living in their own runs).
$rng = $d.Content
$found = $rng.Find.Execute("Ya no es necesario añadir .", $true, $false, $false, $false, $false, $true, 1, $false)

if ($found) {
    # Overwrite the matched range directly (rather than using Find's
    # Replace) so that Word's smart-quote autocorrect doesn't mangle the
    # straight double quotes inside the inserted markup. This also merges
    # the three original runs into the single run the new text lives in.
    $rng.Text = 'Ya no es necesario añadir <span class="caret"></span>.'
}
